# Updated via Streamlit Approval System
# Duplicate the last data row (row 22) into a new row 23, then update the
# fields that differ for this new pending-approval entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the entire last row (22) down to the new row (23) so that all
# shared/common fields (employee, company, bank details, etc.) and cell
# formatting (e.g. the date style on column C) come across intact.
$ws.Range("A22:AO22").Copy($ws.Range("A23"))

# Make sure the date cell keeps its original date/time format (Copy should
# already do this, but set it explicitly so we don't accidentally add a
# new number format to the stylesheet).
$ws.Range("C23").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Now overwrite the cells that are specific to this new record.
$ws.Range("L23").Value2 = "aff807ee-c44c-44d8-ad9e-dd7e229333c7"
$ws.Range("V23").Value2 = 238000
$ws.Range("X23").Value2 = "Payment of ISD RPA_UNIQUE_ID : b12b5ded-ce35-46cd-9876-f7bb854714b9"
$ws.Range("Y23").Value2 = "Construction of New RO (Resitement of M/s Jane Austin,`nKollam to Puthenkurish, Ernakulam and Conversion from B to A at Village) at Thiruvaniyoor, Taluk:`nKunnathunadu, Distt.: Ernakulam, Post- Puthenkurish- 682308, Cochin Divisional Office under Kerala State`nOffice."
$ws.Range("Z23").Value2 = "PAYMENT OF ISD"
$ws.Range("AA23").Value2 = "midhuncraju12@gmail.com"

# The narration text in Y23 contains embedded line breaks; reset the row
# height back to the sheet default (no explicit/custom height) instead of
# leaving the auto row-height bump that a multi-line value would trigger.
$ws.Rows.Item(23).AutoFit()
